$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '66.241.84'
$ws.Range("E2").Value = '  +2.48%  '
$ws.Range("D3").Value = '3.245.23'
$ws.Range("E3").Value = '  +5.20%  '
$ws.Range("D4").Value = '''1.00'
$ws.Range("E4").Value = '  +0.15%  '
$ws.Range("D5").Value = '''574.03'
$ws.Range("E5").Value = '  +1.73%  '
$ws.Range("D6").Value = '''154.69'
$ws.Range("E6").Value = '  +8.19%  '
$ws.Range("D8").Value = '3.234.85'
$ws.Range("E8").Value = '  +5.08%  '
$ws.Range("D9").Value = '''0.515'
$ws.Range("E9").Value = '  +3.96%  '
$ws.Range("D10").Value = '''7.17'
$ws.Range("E10").Value = '  +11.25%  '
$ws.Range("D11").Value = '''0.167'
$ws.Range("E11").Value = '  +4.98%  '
$ws.Range("D12").Value = '''0.485'
$ws.Range("E12").Value = '  +3.94%  '
$ws.Range("D13").Value = '''38.21'
$ws.Range("E13").Value = '  +6.71%  '
$ws.Range("D14").Value = '''0.0000236'
$ws.Range("E14").Value = '  +3.89%  '
$ws.Range("D15").Value = '3.759.96'
$ws.Range("E15").Value = '  +5.31%  '
$ws.Range("D16").Value = '66.276.71'
$ws.Range("E16").Value = '  +2.70%  '
$ws.Range("D17").Value = '''547.49'
$ws.Range("E17").Value = '  +10.09%  '
$ws.Range("D18").Value = '3.243.41'
$ws.Range("E18").Value = '  +5.29%  '
$ws.Range("D19").Value = '''0.115'
$ws.Range("E19").Value = '  +2.89%  '
$ws.Range("D20").Value = '''7.06'
$ws.Range("E20").Value = '  +5.83%  '
$ws.Range("D21").Value = '''14.53'
$ws.Range("E21").Value = '  +5.31%  '
$ws.Range("D22").Value = '''0.739'
$ws.Range("E22").Value = '  +7.18%  '
$ws.Range("D23").Value = '''7.76'
$ws.Range("E23").Value = '  +7.21%  '
$ws.Range("D24").Value = '''13.50'
$ws.Range("E24").Value = '  +6.05%  '
$ws.Range("D25").Value = '''81.80'
$ws.Range("E25").Value = '  +3.95%  '
$ws.Range("E26").Value = '  -0.28%  '
$ws.Range("E27").Value = '  +17.18%  '
$ws.Range("E28").Value = '  +4.42%  '
$ws.Range("D29").Value = '''2.27'
$ws.Range("E29").Value = '  +8.77%  '
$ws.Range("D30").Value = '''28.00'
$ws.Range("E30").Value = '  +5.38%  '
$ws.Range("D31").Value = '''2.80'
$ws.Range("E31").Value = '  +5.93%  '
$ws.Range("D32").Value = '''1.00'
$ws.Range("E32").Value = '  -0.05%  '
$ws.Range("E33").Value = '  +4.63%  '
$ws.Range("D34").Value = '''574.89'
$ws.Range("E34").Value = '  +11.01%  '
$ws.Range("D35").Value = '''5.83'
$ws.Range("E35").Value = '  +5.35%  '
$ws.Range("D36").Value = '''6.46'
$ws.Range("E36").Value = '  +7.72%  '
$ws.Range("D37").Value = '''0.0468'
$ws.Range("E37").Value = '  +15.27%  '
$ws.Range("D38").Value = '''54.83'
$ws.Range("E38").Value = '  +2.95%  '
$ws.Range("D39").Value = '''0.0873'
$ws.Range("E39").Value = '  +9.09%  '
$ws.Range("D40").Value = '''3.07'
$ws.Range("E40").Value = '  +15.12%  '
$ws.Range("E41").Value = '  +4.66%  '
$ws.Range("D42").Value = '3.145.26'
$ws.Range("E42").Value = '  +6.51%  '
$ws.Range("D43").Value = '''8.64'
$ws.Range("E43").Value = '  +3.02%  '
$ws.Range("D44").Value = '''2.38'
$ws.Range("E44").Value = '  +12.14%  '
$ws.Range("D45").Value = '''0.275'
$ws.Range("E45").Value = '  +11.71%  '
$ws.Range("D46").Value = '''26.94'
$ws.Range("E46").Value = '  +7.18%  '
$ws.Range("B47").Value = 'PEPE'
$ws.Range("C47").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D47").Value = '0.0₃0563'
$ws.Range("E47").Value = '  +3.55%  '
$ws.Range("B48").Value = 'USDe'
$ws.Range("C48").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D48").Value = '''1.00'
$ws.Range("E48").Value = '  +0.10%  '
$ws.Range("E49").Value = '  +4.61%  '
$ws.Range("B50").Value = 'ThetaToken'
$ws.Range("C50").Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range("D50").Value = '''2.26'
$ws.Range("E50").Value = '  +8.74%  '
$ws.Range("B51").Value = 'Monero'
$ws.Range("C51").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D51").Value = '''122.45'
$ws.Range("E51").Value = '  +0.54%  '
